$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: D2 changes from "ECs" to "MuSCs"; recalculated numeric values
$ws.Range("D2").Value = "MuSCs"
$ws.Range("I2").Value = 0.9300694554254023
$ws.Range("J2").Value = 0.9300694554254023
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4426103333333333
$ws.Range("N2").Value = 1.327831
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1.109899114220444
$ws.Range("R2").Value = 9.989092027984
$ws.Range("S2").Value = 0.9300694554254023
$ws.Range("T2").Value = 0.9300694554254023

# Update row 3: A3 changes from "FAPs" to "MuSCs"; recalculated numeric values
$ws.Range("A3").Value = "MuSCs"
$ws.Range("G3").Value = 0.1885443333333333
$ws.Range("H3").Value = 0.5656329999999999
$ws.Range("I3").Value = 0.06993054457459773
$ws.Range("J3").Value = 0.06993054457459771
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.4426103333333333
$ws.Range("N3").Value = 1.327831
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.08345167022477777
$ws.Range("R3").Value = 0.7510650320229999
$ws.Range("S3").Value = 0.06993054457459773
$ws.Range("T3").Value = 0.06993054457459771

# Delete rows 4 and 5 entirely
$ws.Range("A4:T5").EntireRow.Delete()
